# Update the dSF column (F) values for the gausman_kevin 2024 save_data sheet.
# Row -> new F value mapping, derived from the target diff. Rows not listed
# (8, 14, 26) are left unchanged because their F value already matched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 3
    3  = 2
    4  = 1
    5  = -2
    6  = -2
    7  = 3
    9  = 5
    10 = 5
    11 = 2
    12 = 4
    13 = 1
    15 = -1
    16 = -3
    17 = 5
    18 = -1
    19 = -3
    20 = 7
    21 = -4
    22 = 3
    23 = 1
    24 = 1
    25 = -6
    27 = 2
    28 = -1
    29 = 1
    30 = -2
    31 = -5
    32 = 6
    33 = 1
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 6).Value = $values[$row]
}
